# Commit: "Removed Blockchain Reference in Chain."
#
# 1) On the two "List<BlockLocation>" labels (slide 1 and slide 2) drop the
#    stray trailing <a:endParaRPr> run (PowerPoint leaves an explicit
#    end-of-paragraph run-properties element after the only run is
#    deleted/retyped; remove it by deleting the paragraph text and
#    re-inserting it, which collapses back down to just the <a:r> run).
# 2) Rename the "Chain" label (nested two levels deep inside groups on
#    slide 2) to "MainChain".

$p = $ppt.ActivePresentation

function Remove-TrailingEndParaRPr($shape) {
    $tr = $shape.TextFrame.TextRange
    $text = $tr.Text
    $tr.Delete()
    $shape.TextFrame.TextRange.InsertAfter($text) | Out-Null
}

# --- Slide 1: "List<BlockLocation>" textbox (TextBox 41) ---
$slide1 = $p.Slides.Item(1)
$blockLoc1 = $slide1.Shapes.Item("TextBox 41")
Remove-TrailingEndParaRPr $blockLoc1

# --- Slide 2: "List<BlockLocation>" textbox (TextBox 21) ---
$slide2 = $p.Slides.Item(2)
$blockLoc2 = $slide2.Shapes.Item("TextBox 21")
Remove-TrailingEndParaRPr $blockLoc2

# --- Slide 2: "Chain" -> "MainChain" (Group 38 > GroupItems("TextBox 13")) ---
$chainGroup = $slide2.Shapes.Item("Group 38")
$chainLabel = $chainGroup.GroupItems.Item("TextBox 13")
$chainLabel.TextFrame.TextRange.Text = "MainChain"
